$wb = $excel.ActiveWorkbook

# --- Add Natick sheet (after Hudson, the last existing sheet) ---
$hudson = $wb.Worksheets.Item($wb.Worksheets.Count)
$newNatick = $wb.Worksheets.Add([System.Type]::Missing, $hudson)
$newNatick.Name = "Natick"

$newNatick.Range("A1").Value = "Day"
$newNatick.Range("B1").Value = "Date"
$newNatick.Range("C1").Value = "Time"
$newNatick.Range("D1").Value = "Location"
$newNatick.Range("E1").Value = "Address"
$newNatick.Range("A2").Value = "Monday"
$newNatick.Range("B2").Value = 43395
$newNatick.Range("B2").NumberFormat = "d-mmm"
$newNatick.Range("C2").Value = "8:00 AM - 5:00 PM"
$newNatick.Range("D2").Value = "NATICK TOWN HALL"
$newNatick.Range("E2").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A3").Value = "Tuesday"
$newNatick.Range("B3").Value = 43396
$newNatick.Range("B3").NumberFormat = "d-mmm"
$newNatick.Range("C3").Value = "8:00 AM - 5:00 PM"
$newNatick.Range("D3").Value = "NATICK TOWN HALL"
$newNatick.Range("E3").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A4").Value = "Wednesday"
$newNatick.Range("B4").Value = 43397
$newNatick.Range("B4").NumberFormat = "d-mmm"
$newNatick.Range("C4").Value = "8:00 AM - 5:00 PM"
$newNatick.Range("D4").Value = "NATICK TOWN HALL"
$newNatick.Range("E4").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A5").Value = "Thursday"
$newNatick.Range("B5").Value = 43398
$newNatick.Range("B5").NumberFormat = "d-mmm"
$newNatick.Range("C5").Value = "8:00 AM - 7:00 PM"
$newNatick.Range("D5").Value = "NATICK TOWN HALL"
$newNatick.Range("E5").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A6").Value = "Friday"
$newNatick.Range("B6").Value = 43399
$newNatick.Range("B6").NumberFormat = "d-mmm"
$newNatick.Range("C6").Value = "8:00 AM - 4:00 PM"
$newNatick.Range("D6").Value = "NATICK TOWN HALL"
$newNatick.Range("E6").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A7").Value = "Saturday"
$newNatick.Range("B7").Value = 43400
$newNatick.Range("B7").NumberFormat = "d-mmm"
$newNatick.Range("C7").Value = "9:00 AM - 3:00 PM"
$newNatick.Range("D7").Value = "NATICK TOWN HALL"
$newNatick.Range("E7").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A8").Value = "Monday"
$newNatick.Range("B8").Value = 43402
$newNatick.Range("B8").NumberFormat = "d-mmm"
$newNatick.Range("C8").Value = "8:00 AM - 5:00 PM"
$newNatick.Range("D8").Value = "NATICK TOWN HALL"
$newNatick.Range("E8").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A9").Value = "Tuesday"
$newNatick.Range("B9").Value = 43403
$newNatick.Range("B9").NumberFormat = "d-mmm"
$newNatick.Range("C9").Value = "8:00 AM - 5:00 PM"
$newNatick.Range("D9").Value = "NATICK TOWN HALL"
$newNatick.Range("E9").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A10").Value = "Wednesday"
$newNatick.Range("B10").Value = 43404
$newNatick.Range("B10").NumberFormat = "d-mmm"
$newNatick.Range("C10").Value = "8:00 AM - 5:00 PM"
$newNatick.Range("D10").Value = "NATICK TOWN HALL"
$newNatick.Range("E10").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A11").Value = "Thursday"
$newNatick.Range("B11").Value = 43405
$newNatick.Range("B11").NumberFormat = "d-mmm"
$newNatick.Range("C11").Value = "8:00 AM - 7:00 PM"
$newNatick.Range("D11").Value = "NATICK TOWN HALL"
$newNatick.Range("E11").Value = "13 EAST CENTRAL ST."
$newNatick.Range("A12").Value = "Friday"
$newNatick.Range("B12").Value = 43406
$newNatick.Range("B12").NumberFormat = "d-mmm"
$newNatick.Range("C12").Value = "8:00 AM - 4:00 PM"
$newNatick.Range("D12").Value = "NATICK TOWN HALL"
$newNatick.Range("E12").Value = "13 EAST CENTRAL ST."

$newNatick.Columns.Item(3).ColumnWidth = 24.7
$newNatick.Range("A1:E12").Select()

# --- Add Arlington sheet (after Natick) ---
$newArlington = $wb.Worksheets.Add([System.Type]::Missing, $newNatick)
$newArlington.Name = "Arlington"

$newArlington.Range("A1").Value = "Day"
$newArlington.Range("B1").Value = "Date"
$newArlington.Range("C1").Value = "Time"
$newArlington.Range("D1").Value = "Location"
$newArlington.Range("E1").Value = "Address"
$newArlington.Range("A2").Value = "Monday"
$newArlington.Range("B2").Value = 43395
$newArlington.Range("B2").NumberFormat = "d-mmm"
$newArlington.Range("C2").Value = "8:00 AM - 4:00 PM"
$newArlington.Range("D2").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E2").Value = "730 MASSACHUSETTS AVENUE"
$newArlington.Range("A3").Value = "Tuesday"
$newArlington.Range("B3").Value = 43396
$newArlington.Range("B3").NumberFormat = "d-mmm"
$newArlington.Range("C3").Value = "8:00 AM - 4:00 PM"
$newArlington.Range("D3").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E3").Value = "730 MASSACHUSETTS AVENUE"
$newArlington.Range("A4").Value = "Wednesday"
$newArlington.Range("B4").Value = 43397
$newArlington.Range("B4").NumberFormat = "d-mmm"
$newArlington.Range("C4").Value = "8:00 AM - 4:00 PM"
$newArlington.Range("D4").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E4").Value = "730 MASSACHUSETTS AVENUE"
$newArlington.Range("A5").Value = "Thursday"
$newArlington.Range("B5").Value = 43398
$newArlington.Range("B5").NumberFormat = "d-mmm"
$newArlington.Range("C5").Value = "8:00 AM - 7:00 PM"
$newArlington.Range("D5").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E5").Value = "730 MASSACHUSETTS AVENUE"
$newArlington.Range("A6").Value = "Friday"
$newArlington.Range("B6").Value = 43399
$newArlington.Range("B6").NumberFormat = "d-mmm"
$newArlington.Range("C6").Value = "8:00 AM - 12:00 PM"
$newArlington.Range("D6").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E6").Value = "730 MASSACHUSETTS AVENUE"
$newArlington.Range("A7").Value = "Monday"
$newArlington.Range("B7").Value = 43402
$newArlington.Range("B7").NumberFormat = "d-mmm"
$newArlington.Range("C7").Value = "8:00 AM - 4:00 PM"
$newArlington.Range("D7").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E7").Value = "730 MASSACHUSETTS AVENUE"
$newArlington.Range("A8").Value = "Tuesday"
$newArlington.Range("B8").Value = 43403
$newArlington.Range("B8").NumberFormat = "d-mmm"
$newArlington.Range("C8").Value = "8:00 AM - 4:00 PM"
$newArlington.Range("D8").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E8").Value = "730 MASSACHUSETTS AVENUE"
$newArlington.Range("A9").Value = "Wednesday"
$newArlington.Range("B9").Value = 43404
$newArlington.Range("B9").NumberFormat = "d-mmm"
$newArlington.Range("C9").Value = "8:00 AM - 4:00 PM"
$newArlington.Range("D9").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E9").Value = "730 MASSACHUSETTS AVENUE"
$newArlington.Range("A10").Value = "Thursday"
$newArlington.Range("B10").Value = 43405
$newArlington.Range("B10").NumberFormat = "d-mmm"
$newArlington.Range("C10").Value = "8:00 AM - 7:00 PM"
$newArlington.Range("D10").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E10").Value = "730 MASSACHUSETTS AVENUE"
$newArlington.Range("A11").Value = "Friday"
$newArlington.Range("B11").Value = 43406
$newArlington.Range("B11").NumberFormat = "d-mmm"
$newArlington.Range("C11").Value = "8:00 AM - 12:00 PM"
$newArlington.Range("D11").Value = "TOWN HALL AUDITORIUM"
$newArlington.Range("E11").Value = "730 MASSACHUSETTS AVENUE"

$newArlington.Range("A1:E11").Select()

# --- Add Southborough sheet (after Arlington) ---
$newSouthborough = $wb.Worksheets.Add([System.Type]::Missing, $newArlington)
$newSouthborough.Name = "Southborough"

$newSouthborough.Range("A1").Value = "Day"
$newSouthborough.Range("B1").Value = "Date"
$newSouthborough.Range("C1").Value = "Time"
$newSouthborough.Range("D1").Value = "Location"
$newSouthborough.Range("E1").Value = "Address"
$newSouthborough.Range("A2").Value = "Monday"
$newSouthborough.Range("B2").Value = 43395
$newSouthborough.Range("B2").NumberFormat = "d-mmm"
$newSouthborough.Range("C2").Value = "8:00 AM - 5:00 PM"
$newSouthborough.Range("D2").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E2").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A3").Value = "Tuesday"
$newSouthborough.Range("B3").Value = 43396
$newSouthborough.Range("B3").NumberFormat = "d-mmm"
$newSouthborough.Range("C3").Value = "8:00 AM - 7:00 PM"
$newSouthborough.Range("D3").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E3").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A4").Value = "Wednesday"
$newSouthborough.Range("B4").Value = 43397
$newSouthborough.Range("B4").NumberFormat = "d-mmm"
$newSouthborough.Range("C4").Value = "8:00 AM - 5:00 PM"
$newSouthborough.Range("D4").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E4").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A5").Value = "Thursday"
$newSouthborough.Range("B5").Value = 43398
$newSouthborough.Range("B5").NumberFormat = "d-mmm"
$newSouthborough.Range("C5").Value = "8:00 AM - 5:00 PM"
$newSouthborough.Range("D5").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E5").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A6").Value = "Friday"
$newSouthborough.Range("B6").Value = 43399
$newSouthborough.Range("B6").NumberFormat = "d-mmm"
$newSouthborough.Range("C6").Value = "8:00 AM - 5:00 PM"
$newSouthborough.Range("D6").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E6").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A7").Value = "Saturday"
$newSouthborough.Range("B7").Value = 43400
$newSouthborough.Range("B7").NumberFormat = "d-mmm"
$newSouthborough.Range("C7").Value = "8:00 AM - 6:00 PM"
$newSouthborough.Range("D7").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E7").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A8").Value = "Monday"
$newSouthborough.Range("B8").Value = 43402
$newSouthborough.Range("B8").NumberFormat = "d-mmm"
$newSouthborough.Range("C8").Value = "8:00 AM - 5:00 PM"
$newSouthborough.Range("D8").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E8").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A9").Value = "Tuesday"
$newSouthborough.Range("B9").Value = 43403
$newSouthborough.Range("B9").NumberFormat = "d-mmm"
$newSouthborough.Range("C9").Value = "8:00 AM - 7:00 PM"
$newSouthborough.Range("D9").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E9").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A10").Value = "Wednesday"
$newSouthborough.Range("B10").Value = 43404
$newSouthborough.Range("B10").NumberFormat = "d-mmm"
$newSouthborough.Range("C10").Value = "8:00 AM - 5:00 PM"
$newSouthborough.Range("D10").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E10").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A11").Value = "Thursday"
$newSouthborough.Range("B11").Value = 43405
$newSouthborough.Range("B11").NumberFormat = "d-mmm"
$newSouthborough.Range("C11").Value = "8:00 AM - 5:00 PM"
$newSouthborough.Range("D11").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E11").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"
$newSouthborough.Range("A12").Value = "Friday"
$newSouthborough.Range("B12").Value = 43406
$newSouthborough.Range("B12").NumberFormat = "d-mmm"
$newSouthborough.Range("C12").Value = "8:00 AM - 6:00 PM"
$newSouthborough.Range("D12").Value = "SOUTHBOROUGH TOWN HOUSE"
$newSouthborough.Range("E12").Value = "17 COMMON STREET SOUTHBOROUGH, MA 01772"

$newSouthborough.Range("J18").Select()
$newSouthborough.Activate()

Write-Output "Added Natick, Arlington, Southborough sheets"
